$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column G ("K") values for rows 3-8 per the diff
$ws.Range("G3").Value = 9
$ws.Range("G4").Value = 7
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 7
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 3
